# Insert a new data row at row 92 (pushing existing rows 92-198 down to 93-199)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(92).Insert()

$ws.Cells.Item(92, 1).Value  = 1
$ws.Cells.Item(92, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(92, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(92, 4).Value  = 44483
$ws.Cells.Item(92, 5).Value  = 15
$ws.Cells.Item(92, 6).Value  = 100114013
$ws.Cells.Item(92, 7).Value  = "Zanahoria"
$ws.Cells.Item(92, 8).Value  = "Sin especificar"
$ws.Cells.Item(92, 9).Value  = "Primera"
$ws.Cells.Item(92, 10).Value = 50
$ws.Cells.Item(92, 11).Value = 9000
$ws.Cells.Item(92, 12).Value = 10000
$ws.Cells.Item(92, 13).Value = 9500
$ws.Cells.Item(92, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(92, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(92, 16).Value = 380
$ws.Cells.Item(92, 17).Value = 25
$ws.Cells.Item(92, 18).Value = "Hortaliza"
